$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3266
$ws1.Range("F3").Value = 9
$ws1.Range("F4").Value = 56
$ws1.Range("F5").Value = 1202
$ws1.Range("F6").Value = 311

# Sheet "全部类型" (All types) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3266
$ws4.Range("F3").Value = 9
$ws4.Range("F4").Value = 56
$ws4.Range("F5").Value = 1202
$ws4.Range("F7").Value = 311
